$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.161.37"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.549.78"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'619.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").Value = "'174.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "3.543.04"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "'7.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "'0.592"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "'46.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "4.119.71"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "'8.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "'616.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "71.182.37"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.540.69"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("D22").Value = "'0.893"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").Value = "'9.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.23%  "
$ws.Range("D24").Value = "'15.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'98.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").Value = "'3.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "'33.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "'3.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").Value = "'8.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'6.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "'624.66"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.75%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").Value = "'10.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "'3.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("D39").Value = "'0.0480"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").Value = "'57.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("D44").Value = "3.384.71"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "'32.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("D50").Value = "'133.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "
